# Refresh the "cryptos" price/volume snapshot (GitHub Actions style update).
# Column D ("Price") cells are stored as text in this sheet (values like
# "30.365.24" or "1.0000" aren't valid numbers / would lose their exact
# formatting if Excel auto-typed them), so each D assignment is entered
# with a leading apostrophe to force text entry, then the cell style is
# reset to "Normal" so no stray text-format style is left behind.
# Column E ("Volume(1h)") cells are plain padded percentage strings and
# assign safely as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.369.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "'1.871.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.4703"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.21%  "
$ws.Range("D8").Value = "'0.2876"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'22.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").Value = "'0.07771"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "'1.869.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "'96.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "'0.7246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").Value = "'5.123"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "'279.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").Value = "'30.354.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "'12.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").Value = "'1.0000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "'0.000007491"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "'2.111.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'5.227"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "'6.229"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").Value = "'163.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").Value = "'9.054"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("D27").Value = "'18.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").Value = "'1.874"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("D29").Value = "'1.322"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").Value = "'4.210"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("D33").Value = "'4.113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "'0.04809"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "'1.119"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").Value = "'0.6893"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'2.716"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'0.01874"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "'2.807"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").Value = "'6.222"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("D41").Value = "'74.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").Value = "'0.4230"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "'1.932"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "'0.9993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "'0.8291"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").Value = "'100.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "'9.633"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("D48").Value = "'35.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'6.947"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "'903.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "'0.05720"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.78%  "
